$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Omg"
$ws.Range("C2").Value = "Rtn4rl1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.084811666666667
$ws.Range("H2").Value = 15.254435
$ws.Range("I2").Value = 0.5108142255059609
$ws.Range("J2").Value = 0.5108142255059609
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.128483
$ws.Range("N2").Value = 0.385449
$ws.Range("O2").Value = 0.01580103135560779
$ws.Range("P2").Value = 0.01580103135560779
$ws.Range("Q2").Value = 0.6533118573683333
$ws.Range("R2").Value = 5.879806716315001
$ws.Range("S2").Value = 0.008071391594110197
$ws.Range("T2").Value = 0.008071391594110197

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Omg"
$ws.Range("C3").Value = "Rtn4rl1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.084811666666667
$ws.Range("H3").Value = 15.254435
$ws.Range("I3").Value = 0.5108142255059609
$ws.Range("J3").Value = 0.5108142255059609
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.550946666666667
$ws.Range("N3").Value = 13.65284
$ws.Range("O3").Value = 0.5596822223772701
$ws.Range("P3").Value = 0.55968222237727
$ws.Range("Q3").Value = 23.14070670504445
$ws.Range("R3").Value = 208.2663603454
$ws.Range("S3").Value = 0.2858936409531002
$ws.Range("T3").Value = 0.2858936409531001

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Omg"
$ws.Range("C4").Value = "Rtn4rl1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.084811666666667
$ws.Range("H4").Value = 15.254435
$ws.Range("I4").Value = 0.5108142255059609
$ws.Range("J4").Value = 0.5108142255059609
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.451875
$ws.Range("N4").Value = 10.355625
$ws.Range("O4").Value = 0.4245167462671222
$ws.Range("P4").Value = 0.4245167462671222
$ws.Range("Q4").Value = 17.552134271875
$ws.Range("R4").Value = 157.969208446875
$ws.Range("S4").Value = 0.2168491929587505
$ws.Range("T4").Value = 0.2168491929587505

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Omg"
$ws.Range("C5").Value = "Rtn4rl1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.689016666666667
$ws.Range("H5").Value = 8.06705
$ws.Range("I5").Value = 0.2701354653822224
$ws.Range("J5").Value = 0.2701354653822224
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.128483
$ws.Range("N5").Value = 0.385449
$ws.Range("O5").Value = 0.01580103135560779
$ws.Range("P5").Value = 0.01580103135560779
$ws.Range("Q5").Value = 0.3454929283833333
$ws.Range("R5").Value = 3.10943635545
$ws.Range("S5").Value = 0.0042684189587662
$ws.Range("T5").Value = 0.0042684189587662

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Omg"
$ws.Range("C6").Value = "Rtn4rl1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.689016666666667
$ws.Range("H6").Value = 8.06705
$ws.Range("I6").Value = 0.2701354653822224
$ws.Range("J6").Value = 0.2701354653822224
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.550946666666667
$ws.Range("N6").Value = 13.65284
$ws.Range("O6").Value = 0.5596822223772701
$ws.Range("P6").Value = 0.55968222237727
$ws.Range("Q6").Value = 12.23757143577778
$ws.Range("R6").Value = 110.138142922
$ws.Range("S6").Value = 0.1511900176080404
$ws.Range("T6").Value = 0.1511900176080403

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Omg"
$ws.Range("C7").Value = "Rtn4rl1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.689016666666667
$ws.Range("H7").Value = 8.06705
$ws.Range("I7").Value = 0.2701354653822224
$ws.Range("J7").Value = 0.2701354653822224
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.451875
$ws.Range("N7").Value = 10.355625
$ws.Range("O7").Value = 0.4245167462671222
$ws.Range("P7").Value = 0.4245167462671222
$ws.Range("Q7").Value = 9.28214940625
$ws.Range("R7").Value = 83.53934465625
$ws.Range("S7").Value = 0.1146770288154159
$ws.Range("T7").Value = 0.1146770288154159

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Omg"
$ws.Range("C8").Value = "Rtn4rl1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.180498333333333
$ws.Range("H8").Value = 6.541495
$ws.Range("I8").Value = 0.2190503091118167
$ws.Range("J8").Value = 0.2190503091118167
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.128483
$ws.Range("N8").Value = 0.385449
$ws.Range("O8").Value = 0.01580103135560779
$ws.Range("P8").Value = 0.01580103135560779
$ws.Range("Q8").Value = 0.2801569673616667
$ws.Range("R8").Value = 2.521412706255
$ws.Range("S8").Value = 0.003461220802731395
$ws.Range("T8").Value = 0.003461220802731395

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Omg"
$ws.Range("C9").Value = "Rtn4rl1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.180498333333333
$ws.Range("H9").Value = 6.541495
$ws.Range("I9").Value = 0.2190503091118167
$ws.Range("J9").Value = 0.2190503091118167
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.550946666666667
$ws.Range("N9").Value = 13.65284
$ws.Range("O9").Value = 0.5596822223772701
$ws.Range("P9").Value = 0.55968222237727
$ws.Range("Q9").Value = 9.923331621755558
$ws.Range("R9").Value = 89.30998459580002
$ws.Range("S9").Value = 0.1225985638161295
$ws.Range("T9").Value = 0.1225985638161295

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Omg"
$ws.Range("C10").Value = "Rtn4rl1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.180498333333333
$ws.Range("H10").Value = 6.541495
$ws.Range("I10").Value = 0.2190503091118167
$ws.Range("J10").Value = 0.2190503091118167
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.451875
$ws.Range("N10").Value = 10.355625
$ws.Range("O10").Value = 0.4245167462671222
$ws.Range("P10").Value = 0.4245167462671222
$ws.Range("Q10").Value = 7.526807684375
$ws.Range("R10").Value = 67.74126915937501
$ws.Range("S10").Value = 0.09299052449295578
$ws.Range("T10").Value = 0.09299052449295576
